$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing row 4 timestamp (tiny precision correction from the
# scheduled task re-running at save time).
$ws.Range("A4").Value2 = 45865.12527226852

# Append the new row of sensor readings captured by the scheduled task.
$ws.Range("A5").Value2 = 45865.16689058254
$ws.Range("A5").NumberFormat = $ws.Range("A4").NumberFormat

$ws.Range("B5").Value2 = 2025
$ws.Range("C5").Value2 = 30
$ws.Range("D5").Value2 = 13.44
$ws.Range("E5").Value2 = 90.88
$ws.Range("F5").Value2 = 0
$ws.Range("G5").Value2 = 4.02
$ws.Range("H5").Value2 = "WNW"
$ws.Range("I5").Value2 = 0
$ws.Range("J5").Value2 = "04:00:19"
